$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text would otherwise be auto-parsed as a number by Excel;
# force Text format first so the scraped string (e.g. "589.27") round-trips as text,
# matching the sheets existing inline-string "Price" column convention.
$textPriceCells = @(
    "D5", "D6", "D8", "D10", "D12", "D13", "D14", "D17", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D42", "D44", "D45", "D46", "D49", "D50", "D51"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update price (D) and volume-change (E) columns for rows with new scraped values
$ws.Range("D2").Value = "66.403.24"
$ws.Range("E2").Value = "  -4.00%  "
$ws.Range("D3").Value = "3.573.04"
$ws.Range("E3").Value = "  -4.35%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "589.27"
$ws.Range("E5").Value = "  -4.73%  "
$ws.Range("D6").Value = "186.51"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "3.564.84"
$ws.Range("E7").Value = "  -4.57%  "
$ws.Range("D8").Value = "0.614"
$ws.Range("E8").Value = "  -4.09%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "0.671"
$ws.Range("E11").Value = "  -9.58%  "
$ws.Range("D12").Value = "53.65"
$ws.Range("E12").Value = "  -5.79%  "
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  -12.09%  "
$ws.Range("D14").Value = "9.77"
$ws.Range("E14").Value = "  -8.42%  "
$ws.Range("D15").Value = "4.140.68"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").Value = "3.570.14"
$ws.Range("E16").Value = "  -4.28%  "
$ws.Range("D17").Value = "0.125"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "12.24"
$ws.Range("E18").Value = "  -6.53%  "
$ws.Range("D19").Value = "18.29"
$ws.Range("E19").Value = "  -5.95%  "
$ws.Range("D20").Value = "66.319.20"
$ws.Range("E20").Value = "  -3.83%  "
$ws.Range("E21").Value = "  -7.09%  "
$ws.Range("D22").Value = "397.15"
$ws.Range("D23").Value = "4.38"
$ws.Range("E23").Value = "  -6.58%  "
$ws.Range("D24").Value = "85.63"
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("D25").Value = "11.26"
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("D26").Value = "2.89"
$ws.Range("E26").Value = "  -5.77%  "
$ws.Range("D27").Value = "12.49"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  -6.86%  "
$ws.Range("D30").Value = "8.98"
$ws.Range("E30").Value = "  -7.18%  "
$ws.Range("D31").Value = "31.07"
$ws.Range("E31").Value = "  -6.71%  "
$ws.Range("D32").Value = "7.09"
$ws.Range("E32").Value = "  -3.54%  "
$ws.Range("D33").Value = "12.14"
$ws.Range("E33").Value = "  -4.80%  "
$ws.Range("D34").Value = "617.91"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").Value = "63.54"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").Value = "0.113"
$ws.Range("D37").Value = "41.43"
$ws.Range("E37").Value = "  -7.07%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "0.394"
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("E40").Value = "  -13.37%  "
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "3.023.26"
$ws.Range("E43").Value = "  +6.54%  "
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -8.37%  "
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  -5.19%  "
$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").Value = "  -8.48%  "
$ws.Range("E47").Value = "  -7.06%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "8.60"
$ws.Range("E49").Value = "  -7.18%  "

# Rows 50 and 51 swapped ranking order (Monero now ranked above Stacks) with updated values
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "138.68"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "2.77"
$ws.Range("E51").Value = "  -0.36%  "
